# Apply the changes described by the diff to the "Dados" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dados")

# --- Row 3 updates ---
# C3: "20/01/2023" -> "03/01/2022"  (keep as text, not an Excel date)
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "03/01/2022"

# D3: "00:00" -> "10:00" (keep as text, not an Excel time)
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "10:00"

# E3: 0 -> 50
$ws.Range("E3").Value = 50

# J3: new empty cell appears (no content, just present)
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "x"
$ws.Range("J3").ClearContents()

# --- Row 6 updates ---
# J6 empty placeholder cell is removed
$ws.Range("J6").ClearContents()

# --- Row 7 (new row) ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Jean"

$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "02/01/2022"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "22:45"

$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("K7").Value = "Pendente"

Write-Output "Edits applied."
